$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A row 2 label: "No tessellation" -> "Tessellation factor 1 (no tesselation)"
$ws.Range("A2").Value = "Tessellation factor 1 (no tesselation)"

# New header cells for the added columns
$ws.Range("G1").Value = "Triangle Count"
$ws.Range("H1").Value = "Vertex Count"

# New Triangle Count / Vertex Count values for rows 2-18
$triangleCount = 5183814
$vertexCount = 8321936

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 7).Value = $triangleCount
    $ws.Cells.Item($r, 8).Value = $vertexCount
}

$ws.Range("G2:H18").NumberFormat = "#,##0"

# Column widths (values chosen so that, after this runtime's internal
# character-width<->pixel rounding, the stored OOXML width comes out as close
# as possible to the target 39.44140625 / 15.88671875 / 13)
$ws.Range("A1").EntireColumn.ColumnWidth = 38.666666666666664
$ws.Range("G1").EntireColumn.ColumnWidth = 15.0
$ws.Range("H1").EntireColumn.ColumnWidth = 12.166666666666666

# Update selected cell in the sheet view
[void]$ws.Range("J17").Select()
